# Apply the changes described by the commit:
# "Set up input data for tests and added applied demand test."
#
# This touches the second worksheet ("Contingency Conservation Test"):
#   - B2 changes from a literal value (100) to a formula (=1000000-5)
#     which recalculates the dependent formula in B6 (=B2*B3).
#   - The active selection on that sheet moves from B5 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contingency Conservation Test")

# Make sure this is the active sheet (it was already the tab-selected one).
$ws.Activate()

# B2: replace the literal 100 with a formula "=1000000-5" (evaluates to 999995).
# B6 (=B2*B3) will automatically recalculate to 99999.5.
$ws.Range("B2").Formula = "=1000000-5"

# Move the selected/active cell on the sheet from B5 to B3.
$ws.Range("B3").Select()
